$d = $word.ActiveDocument

# Locate the "Fluxos do Caso de Uso Listar Ocorrências" Heading 1 paragraph
# (there are other, unrelated "Listar "/"Cadastrar" runs elsewhere in the
# document, e.g. in bullet lists and the cover title, so match on the
# paragraph's own text plus its style to stay precise).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Fluxos do Caso de Uso Listar Ocorrências" -and $p.Style.NameLocal -eq "Heading 1") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $rng.Find.ClearFormatting()
    $rng.Find.Execute("Listar ", $true, $false, $false, $false, $false, $true, 1, $false, "Cadastrar ", 2)
}
